$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "x"
$ws.Range("B1").Value = "y"
$ws.Range("C1").Value = "comment"
$ws.Range("D1").Value = "comment_day"
$ws.Range("E1").Value = "comment_floor_id"
$ws.Range("F1").Value = "comment_tag"

$headerRange = $ws.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Data row
$ws.Range("A2").Value = 241
$ws.Range("B2").Value = 93
$ws.Range("C2").Value = "456"
$ws.Range("D2").Value = (Get-Date -Year 2023 -Month 4 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = "241-93-text"
